$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws_0 = $wb.Worksheets.Item("展览")
$ws_0.Range("F2").Value = 2364
$ws_0.Range("F3").Value = 547
$ws_0.Range("F5").Value = 357
$ws_0.Range("F6").Value = 357
$ws_0.Range("F7").Value = 588
$ws_0.Range("F9").Value = 800
$ws_0.Range("F11").Value = 831
$ws_0.Range("F14").Value = 402
$ws_0.Range("F17").Value = 21410
$ws_0.Range("G17").Value = "暂时售罄"
$ws_0.Range("F18").Value = 873
$ws_0.Range("F19").Value = 82
$ws_0.Range("F20").Value = 271
$ws_0.Range("F21").Value = 305
$ws_0.Range("F25").Value = 18
$ws_0.Range("F26").Value = 246
$ws_0.Range("F28").Value = 361
$ws_0.Range("F29").Value = 160

# Sheet: 演出
$ws_1 = $wb.Worksheets.Item("演出")
$ws_1.Range("F6").Value = 208
$ws_1.Range("F7").Value = 229
$ws_1.Range("F8").Value = 3444
$ws_1.Range("F10").Value = 107
$ws_1.Range("F16").Value = 3910

# Sheet: 本地生活
$ws_2 = $wb.Worksheets.Item("本地生活")
$ws_2.Range("F2").Value = 270
$ws_2.Range("F3").Value = 117
$ws_2.Range("F4").Value = 633

# Sheet: 全部类型
$ws_3 = $wb.Worksheets.Item("全部类型")
$ws_3.Range("F2").Value = 270
$ws_3.Range("F3").Value = 117
$ws_3.Range("F5").Value = 2364
$ws_3.Range("F6").Value = 633
$ws_3.Range("F7").Value = 548
$ws_3.Range("F9").Value = 357
$ws_3.Range("F10").Value = 357
$ws_3.Range("F11").Value = 588
$ws_3.Range("F16").Value = 208
$ws_3.Range("F18").Value = 800
$ws_3.Range("F20").Value = 831
$ws_3.Range("F23").Value = 402
$ws_3.Range("F26").Value = 21410
$ws_3.Range("G26").Value = "暂时售罄"
$ws_3.Range("F27").Value = 229
$ws_3.Range("F28").Value = 3444
$ws_3.Range("F30").Value = 107
$ws_3.Range("F32").Value = 873
$ws_3.Range("F33").Value = 82
$ws_3.Range("F34").Value = 271
$ws_3.Range("F37").Value = 305
$ws_3.Range("F41").Value = 18
$ws_3.Range("F44").Value = 246
$ws_3.Range("F46").Value = 361
$ws_3.Range("F47").Value = 160
$ws_3.Range("F48").Value = 3910
